$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the "Mes" (month) column from numeric month values to text month
# abbreviations (Ene., Feb., ... Dic.), row by row so that the new shared-string
# table entries are created in the same order as the target workbook
# (Dic., Nov., Oct., Sep., Ago., Jul., Jun., May., Abr., Mar., Feb., Ene.).
$ws.Range("C6").Value = "Dic."
$ws.Range("C7").Value = "Nov."
$ws.Range("C8").Value = "Oct."
$ws.Range("C9").Value = "Sep."
$ws.Range("C10").Value = "Ago."
$ws.Range("C11").Value = "Jul."
$ws.Range("C12").Value = "Jun."
$ws.Range("C13").Value = "May."
$ws.Range("C14").Value = "Abr."
$ws.Range("C15").Value = "Mar."
$ws.Range("C16").Value = "Feb."
$ws.Range("C17").Value = "Ene."
$ws.Range("C18").Value = "Dic."
$ws.Range("C19").Value = "Nov."
$ws.Range("C20").Value = "Oct."
$ws.Range("C21").Value = "Sep."
$ws.Range("C22").Value = "Ago."
$ws.Range("C23").Value = "Jul."
$ws.Range("C24").Value = "Jun."
$ws.Range("C25").Value = "May."
$ws.Range("C26").Value = "Abr."
$ws.Range("C27").Value = "Mar."
$ws.Range("C28").Value = "Feb."
$ws.Range("C29").Value = "Ene."
$ws.Range("C30").Value = "Dic."
$ws.Range("C31").Value = "Nov."
$ws.Range("C32").Value = "Oct."
$ws.Range("C33").Value = "Sep."
$ws.Range("C34").Value = "Ago."
$ws.Range("C35").Value = "Jul."
$ws.Range("C36").Value = "Jun."
$ws.Range("C37").Value = "May."
$ws.Range("C38").Value = "Abr."
$ws.Range("C39").Value = "Mar."
$ws.Range("C40").Value = "Feb."
$ws.Range("C41").Value = "Ene."
$ws.Range("C42").Value = "Dic."
$ws.Range("C43").Value = "Nov."
$ws.Range("C44").Value = "Oct."
$ws.Range("C45").Value = "Sep."
$ws.Range("C46").Value = "Ago."
$ws.Range("C47").Value = "Jul."
$ws.Range("C48").Value = "Jun."
$ws.Range("C49").Value = "May."
$ws.Range("C50").Value = "Abr."
$ws.Range("C51").Value = "Mar."
$ws.Range("C52").Value = "Feb."
$ws.Range("C53").Value = "Ene."
$ws.Range("C54").Value = "Dic."
$ws.Range("C55").Value = "Nov."
$ws.Range("C56").Value = "Oct."
$ws.Range("C57").Value = "Sep."
$ws.Range("C58").Value = "Ago."
$ws.Range("C59").Value = "Jul."
$ws.Range("C60").Value = "Jun."
$ws.Range("C61").Value = "May."
$ws.Range("C62").Value = "Abr."
$ws.Range("C63").Value = "Mar."
$ws.Range("C64").Value = "Feb."
$ws.Range("C65").Value = "Ene."
$ws.Range("C66").Value = "Dic."
$ws.Range("C67").Value = "Nov."
$ws.Range("C68").Value = "Oct."
$ws.Range("C69").Value = "Sep."
$ws.Range("C70").Value = "Ago."
$ws.Range("C71").Value = "Jul."
$ws.Range("C72").Value = "Jun."
$ws.Range("C73").Value = "May."
$ws.Range("C74").Value = "Abr."
$ws.Range("C75").Value = "Mar."
$ws.Range("C76").Value = "Feb."
$ws.Range("C77").Value = "Ene."
